$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert a small fragment of run-level OOXML at a COLLAPSED range
# that sits at the very start of a paragraph. This engine's Range.InsertXML
# only splices new runs inline when the insertion point is a paragraph
# boundary, so every helper call below is made against a paragraph-start
# position.
# ---------------------------------------------------------------------------
function Insert-RunXmlAtParaStart($pos, [string]$innerRunsXml) {
    $rr = $d.Range($pos, $pos)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rr.InsertXML($pkg)
}

# ===========================================================================
# Edit 1: Carrier Strike Group-73 (CSG-73) / USS George Washington
#         -> Carrier Strike Group-72 (CSG-72) / USS Abraham Lincoln
# The replacement text is split across four runs (same rPr) as in the
# target revision, so rebuild that run sequence explicitly.
# ===========================================================================
$csgAnchor = $d.Content
$found = $csgAnchor.Find.Execute("Carrier Strike Group-73 (CSG-73): USS George Washington + 1x ")
if (-not $found) {
    throw "Could not locate the CSG-73 / George Washington sentence"
}
$csgParaStart = $csgAnchor.Paragraphs(1).Range.Start

$csgRun4 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> + 1x </w:t></w:r>'
$csgRun3 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Abraham Lincoln</w:t></w:r>'
$csgRun2 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">): USS </w:t></w:r>'
$csgRun1 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Carrier Strike Group-72 (CSG-72</w:t></w:r>'

# Insert in reverse order: each insert prepends right after the paragraph
# mark, so inserting 4,3,2,1 in that order yields the final 1,2,3,4 order.
Insert-RunXmlAtParaStart $csgParaStart $csgRun4
Insert-RunXmlAtParaStart $csgParaStart $csgRun3
Insert-RunXmlAtParaStart $csgParaStart $csgRun2
Insert-RunXmlAtParaStart $csgParaStart $csgRun1

# Remove the original (now-superseded) run text.
$d.Content.Find.Execute("Carrier Strike Group-73 (CSG-73): USS George Washington + 1x ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ===========================================================================
# Edit 2: "Special Operations Component Command consist of ..."
#         -> "... consists of ..." and drop the grammar-error proof marks
#         that used to flag "consist".
# ===========================================================================
$socParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Special Operations Component Command consist *") {
        $socParaIndex = $i
        break
    }
}
if ($socParaIndex -eq -1) {
    throw "Could not locate the Special Operations Component Command paragraph"
}

$socPara = $d.Paragraphs($socParaIndex)
$socStart = $socPara.Range.Start
$socEnd = $socPara.Range.End - 1   # stop before the paragraph mark

# Wipe the paragraph's content (this also clears out the stray gramStart/
# gramEnd proofErr markers around "consist") and rebuild it as plain runs.
$d.Range($socStart, $socEnd).Text = ""

$socRun3 = '<w:r w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of 4 Special Operations Task Groups based both in Israel and Turkey.</w:t></w:r>'
$socRun2 = '<w:r w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>consists</w:t></w:r>'
$socRun1 = '<w:r w:rsidRPr="008336C6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Special Operations Component Command </w:t></w:r>'

Insert-RunXmlAtParaStart $socStart $socRun3
Insert-RunXmlAtParaStart $socStart $socRun2
Insert-RunXmlAtParaStart $socStart $socRun1

Write-Output "Done."
